$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E hold values that look numeric (e.g. "398.61", "1.00")
# but must remain plain text, matching the source workbook (all D/E cells are
# stored as inlineStr). Force text format before writing, then clear the
# explicit number-format style so no stray style index is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '56.101.94'
Set-TextValue $ws.Range("E2") '  +8.74%  '
Set-TextValue $ws.Range("D3") '3.220.26'
Set-TextValue $ws.Range("E3") '  +3.65%  '
Set-TextValue $ws.Range("E4") '  -0.04%  '
Set-TextValue $ws.Range("D5") '398.61'
Set-TextValue $ws.Range("E5") '  +3.34%  '
Set-TextValue $ws.Range("D6") '110.32'
Set-TextValue $ws.Range("E6") '  +6.08%  '
Set-TextValue $ws.Range("E7") '  +2.77%  '
Set-TextValue $ws.Range("E8") '  -0.07%  '
Set-TextValue $ws.Range("E9") '  +6.41%  '
Set-TextValue $ws.Range("D10") '39.39'
Set-TextValue $ws.Range("E10") '  +5.86%  '
Set-TextValue $ws.Range("D11") '0.0908'
Set-TextValue $ws.Range("E11") '  +5.90%  '
Set-TextValue $ws.Range("E12") '  +2.10%  '
Set-TextValue $ws.Range("D13") '3.726.61'
Set-TextValue $ws.Range("E13") '  +3.54%  '
Set-TextValue $ws.Range("E14") '  +3.66%  '
Set-TextValue $ws.Range("E15") '  +2.71%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D16") '3.220.22'
Set-TextValue $ws.Range("E16") '  +3.78%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D17") '1.06'
Set-TextValue $ws.Range("E17") '  +5.70%  '
Set-TextValue $ws.Range("D18") '10.73'
Set-TextValue $ws.Range("E18") '  -1.37%  '
Set-TextValue $ws.Range("D19") '55.960.72'
Set-TextValue $ws.Range("E19") '  +8.41%  '
Set-TextValue $ws.Range("E20") '  +2.72%  '
Set-TextValue $ws.Range("E21") '  +6.81%  '
Set-TextValue $ws.Range("D22") '13.10'
Set-TextValue $ws.Range("E22") '  +5.09%  '
Set-TextValue $ws.Range("D23") '305.99'
Set-TextValue $ws.Range("E23") '  +14.61%  '
Set-TextValue $ws.Range("D24") '75.41'
Set-TextValue $ws.Range("E26") '  +1.77%  '
Set-TextValue $ws.Range("E27") '  +4.37%  '
Set-TextValue $ws.Range("E28") '  +4.71%  '
Set-TextValue $ws.Range("E29") '  +3.95%  '
Set-TextValue $ws.Range("D30") '1.00'
Set-TextValue $ws.Range("E30") '  -0.01%  '
Set-TextValue $ws.Range("E31") '  +4.07%  '
Set-TextValue $ws.Range("E32") '  +7.74%  '
Set-TextValue $ws.Range("D33") '0.0490'
Set-TextValue $ws.Range("E33") '  +1.45%  '
Set-TextValue $ws.Range("D34") '36.18'
Set-TextValue $ws.Range("E34") '  +2.55%  '
Set-TextValue $ws.Range("D35") '2.11'
Set-TextValue $ws.Range("E35") '  +2.17%  '
Set-TextValue $ws.Range("D36") '51.50'
Set-TextValue $ws.Range("E36") '  +2.91%  '
Set-TextValue $ws.Range("D37") '3.13'
Set-TextValue $ws.Range("E37") '  +24.12%  '
Set-TextValue $ws.Range("D38") '1.00'
Set-TextValue $ws.Range("E38") '  +0.05%  '
Set-TextValue $ws.Range("E39") '  +4.36%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D40") '134.65'
Set-TextValue $ws.Range("E40") '  +4.19%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D41") '1.93'
Set-TextValue $ws.Range("E41") '  +2.42%  '
Set-TextValue $ws.Range("D42") '4.03'
Set-TextValue $ws.Range("E42") '  +6.74%  '
Set-TextValue $ws.Range("D43") '17.27'
Set-TextValue $ws.Range("E43") '  +3.88%  '
Set-TextValue $ws.Range("E44") '  -1.55%  '
Set-TextValue $ws.Range("E45") '  +3.00%  '
Set-TextValue $ws.Range("D46") '22.44'
Set-TextValue $ws.Range("E46") '  +0.98%  '
Set-TextValue $ws.Range("D47") '2.50'
Set-TextValue $ws.Range("E47") '  -1.35%  '
Set-TextValue $ws.Range("E48") '  +45.86%  '
Set-TextValue $ws.Range("E49") '  +1.05%  '
Set-TextValue $ws.Range("D50") '2.142.38'
Set-TextValue $ws.Range("E50") '  +3.40%  '
Set-TextValue $ws.Range("D51") '0.0364'
Set-TextValue $ws.Range("E51") '  +13.47%  '
